$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.566.56'
$ws.Range("E2").Value = '  +5.33%  '
$ws.Range("D3").Value = '3.473.75'
$ws.Range("E3").Value = '  +5.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '187.27'
$ws.Range("E5").Value = '  +5.44%  '
$ws.Range("D6").Value = '551.53'
$ws.Range("E6").Value = '  +4.93%  '
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("D8").Value = '3.460.93'
$ws.Range("E8").Value = '  +5.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.641'
$ws.Range("E10").Value = '  +5.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("E12").Value = '  +12.41%  '
$ws.Range("D13").Value = '0.0000278'
$ws.Range("E13").Value = '  +7.62%  '
$ws.Range("D14").Value = '9.52'
$ws.Range("E14").Value = '  +4.71%  '
$ws.Range("D15").Value = '4.009.89'
$ws.Range("E15").Value = '  +5.14%  '
$ws.Range("D16").Value = '3.459.11'
$ws.Range("E16").Value = '  +5.20%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.744.29'
$ws.Range("E17").Value = '  +5.75%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.49%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.121'
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").Value = '11.83'
$ws.Range("E20").Value = '  +6.75%  '
$ws.Range("E21").Value = '  +5.84%  '
$ws.Range("D22").Value = '409.28'
$ws.Range("E22").Value = '  +9.40%  '
$ws.Range("D23").Value = '12.01'
$ws.Range("E23").Value = '  +8.27%  '
$ws.Range("D24").Value = '3.92'
$ws.Range("E24").Value = '  +3.66%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '84.73'
$ws.Range("E25").Value = '  +5.18%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '4.23'
$ws.Range("E26").Value = '  +8.89%  '
$ws.Range("D27").Value = '2.93'
$ws.Range("E27").Value = '  +9.71%  '
$ws.Range("D28").Value = '6.26'
$ws.Range("E28").Value = '  +2.77%  '
$ws.Range("D29").Value = '11.86'
$ws.Range("E29").Value = '  +4.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.25%  '
$ws.Range("D31").Value = '30.35'
$ws.Range("E31").Value = '  +5.24%  '
$ws.Range("D32").Value = '685.97'
$ws.Range("E32").Value = '  +7.37%  '
$ws.Range("D33").Value = '6.98'
$ws.Range("E33").Value = '  +5.26%  '
$ws.Range("D34").Value = '11.67'
$ws.Range("E34").Value = '  +3.82%  '
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  +5.46%  '
$ws.Range("D36").Value = '59.36'
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("E37").Value = '  +19.27%  '
$ws.Range("D38").Value = '38.87'
$ws.Range("E38").Value = '  +6.12%  '
$ws.Range("D39").Value = '0.406'
$ws.Range("E39").Value = '  +4.21%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = '2.82'
$ws.Range("E41").Value = '  +14.52%  '
$ws.Range("E42").Value = '  +21.95%  '
$ws.Range("E43").Value = '  +7.89%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = '3.058.61'
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("D46").Value = '2.97'
$ws.Range("E46").Value = '  +10.47%  '
$ws.Range("D47").Value = '0.0422'
$ws.Range("E47").Value = '  +5.92%  '
$ws.Range("D48").Value = '3.26'
$ws.Range("E48").Value = '  +8.65%  '
$ws.Range("E49").Value = '  +3.96%  '
$ws.Range("D50").Value = '2.73'
$ws.Range("E50").Value = '  +11.90%  '
$ws.Range("E51").Value = '  +11.50%  '
